$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- Edit 1: wrap "ConsultUK" (title run) in spellStart/spellEnd proofErr markers ---
$f1 = $d.Content
$found1 = $f1.Find.Execute("ConsultUK", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    # Re-seat as a plain Range so InsertXML replaces (rather than appends after) the match.
    $r1 = $d.Range($f1.Start, $f1.End)
    $xml1 = '<w:p ' + $wNs + '><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="48A21732"><w:rPr><w:rFonts w:ascii="Aptos" w:eastAsia="Aptos" w:hAnsi="Aptos" w:cs="Aptos"/><w:b/><w:bCs/><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr><w:t>ConsultUK</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
    $r1.InsertXML($xml1)
}

# --- Edit 2: split "Jayrup-" run into a proofErr-wrapped "Jayrup" run plus a separate "-" run ---
$f2 = $d.Content
$found2 = $f2.Find.Execute("Jayrup-", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    # Re-seat as a plain Range so InsertXML replaces (rather than appends after) the match.
    $r2 = $d.Range($f2.Start, $f2.End)
    $xml2 = '<w:p ' + $wNs + '><w:proofErr w:type="spellStart"/><w:r><w:t>Jayrup</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>-</w:t></w:r></w:p>'
    $r2.InsertXML($xml2)
}

Write-Output "done"
